$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the "Week 14 ... W1-13)." paragraph - our anchor.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("W1-13).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$week14Para = $rng.Paragraphs(1)

# ------------------------------------------------------------------
# Step 1: insert a new paragraph right after the Week 14 paragraph,
#         holding the "Some video recordings ..." text (this
#         duplicates the text that used to live further down; that
#         original paragraph is repurposed into a heading below).
# ------------------------------------------------------------------
$afterWeek14 = $week14Para.Next()
$afterWeek14Range = $afterWeek14.Range
$afterWeek14Range.Collapse(1)            # wdCollapseStart
$afterWeek14Range.InsertParagraphBefore()

$videoParaNew = $week14Para.Next()
$videoRange = $videoParaNew.Range
$videoTextRange = $d.Range($videoRange.Start, $videoRange.End - 1)
$videoTextRange.Text = "Some video recordings from the previous years "
$videoTextRange.Collapse(0)              # wdCollapseEnd
$videoTextRange.InsertAfter("might be uploaded, but content of the video lectures might not be up to date. We strongly encourage students to come to class and use video recordings for revision purposes.")

# ------------------------------------------------------------------
# Step 2: the paragraph that used to contain "Some video recordings ..."
#         (reached by walking the paragraph chain, since its text is
#         now duplicated above) becomes a Heading2 "An important note
#         on CNY classes".
# ------------------------------------------------------------------
$emptyPara1 = $videoParaNew.Next()
$oldVideoPara = $emptyPara1.Next()
$oldVideoRange = $oldVideoPara.Range
$oldVideoTextRange = $d.Range($oldVideoRange.Start, $oldVideoRange.End - 1)
$oldVideoTextRange.Text = "An important note on CNY classes"
$oldVideoPara.Style = "Heading 2"

# ------------------------------------------------------------------
# Step 3: insert a new paragraph right after the empty paragraph that
#         follows the heading, with the CNY public-holiday note.
# ------------------------------------------------------------------
$emptyPara2 = $oldVideoPara.Next()

$emptyPara2Range = $emptyPara2.Range
$emptyPara2Range.Collapse(0)             # wdCollapseEnd
$emptyPara2Range.InsertParagraphAfter()

$weekNotePara = $emptyPara2.Next()
$weekNoteRange = $weekNotePara.Range
$weekNoteTextRange = $d.Range($weekNoteRange.Start, $weekNoteRange.End - 1)
$weekNoteTextRange.Text = "On Week 4, Monday 12"
$weekNoteTextRange.Collapse(0)           # wdCollapseEnd
$weekNoteTextRange.InsertAfter("th")
$weekNoteTextRange.Font.Superscript = $true
$weekNoteTextRange.Collapse(0)
$weekNoteTextRange.InsertAfter(" of February is a public holiday. Lectures on that day will be cancelled and replaced with a pre-recorded lecture. It will be uploaded on eDimension closer to the date.")

Write-Host "Edit complete"
